$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$data = @(
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(1, 4),
    @(1, 3),
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(1, 4),
    @(1, 5),
    @(1, 6),
    @(1, 4),
    @(1, 6),
    @(1, 4),
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(1, 4),
    @(1, 4),
    @(1, 3),
    @(6, 6),
    @(6, 8),
    @(7, 8),
    @(8, 9),
    @(6, 9),
    @(7, 8),
    @(3, 9),
    @(6, 8),
    @(1, 5),
    @(1, 4),
    @(1, 4),
    @(1, 5),
    @(1, 4),
    @(6, 8),
    @(1, 2)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $data[$r][0]
    $ws.Cells.Item($row, 10).Value = $data[$r][1]
}
